$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 111
$ws.Range("B4").Value = 222
$ws.Range("C4").Value = "Dic"
$ws.Range("D4").Value = "Eri"

$ws.Range("F4").Select()
